# BLE firmware project added
# Fill in the "Megfelelője" (UUID) and "Neve" (Name) columns for several
# GATT table rows, and reserve (but leave blank) the same columns for a
# few more rows further down the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36 - IR Temperature Data characteristic
$ws.Range("C36").Value = "12:21:00:00:00:00:00:00:00:00:B0:00:40:51:04:01:AA:00:F0"
$ws.Range("D36").Value = "IR Temperature Data"

# Row 39 - IR Temperature Config characteristic
$ws.Range("C39").Value = "0A:24:00:00:00:00:00:00:00:00:B0:00:40:51:04:02:AA:00:F0"
$ws.Range("D39").Value = "IR Temperature Config"

# Row 41 - IR Temperature Period characteristic
$ws.Range("C41").Value = "0A:26:00:00:00:00:00:00:00:00:B0:00:40:51:04:03:AA:00:F0"
$ws.Range("D41").Value = "IR Temperature Period"

# Row 101 - OAT Image Identify characteristic
$ws.Range("C101").Value = "1C:63:00:00:00:00:00:00:00:00:B0:00:40:51:04:C1:FF:00:F0"
$ws.Range("D101").Value = "OAT Image Identify"

# Row 100 - OAT Service declaration
$ws.Range("D100").Value = "OAT Service"
$ws.Range("C100").Value = "F000FFC0-0451-4000-B000-000000000000"

# Row 105 - OAT Image Block characteristic
$ws.Range("C105").Value = "1C:67:00:00:00:00:00:00:00:00:B0:00:40:51:04:C2:FF:00:F0"
$ws.Range("D105").Value = "OAT Image Block"

# Row 35 - IR Temperature Service declaration
$ws.Range("C35").Value = "F000AA00-0451-4000-B000-000000000000"
$ws.Range("D35").Value = "IR Temperature Service"

# Reserve (blank, text-formatted) C/D cells for rows 107-112 so they match
# the style already used throughout the rest of the table.
$ws.Range("C107:D112").NumberFormat = "@"

# Leave the selection where the author ended up editing.
$ws.Range("B36").Select()
